$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the obsolete "Version" column (was column Q) ---------------
# Copy Q1's format (the merged title bar's closing right border) onto P1
# first, since P1 becomes the new right-most cell of that merge once Q
# is gone.
$ws.Range("Q1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# This single structural delete re-derives dimension, merged cells,
# row spans and the BI4/CE4 anchor cells -> BH4/CD4 for free.
$ws.Columns("Q").Delete()

# --- Header row (row 2): rewrite the 15 labeled columns ----------------
$ws.Range("B2").Value = " ورژن نرم افزار تستر"
$ws.Range("C2").Value = "اتصال تغذیه "
$ws.Range("D2").Value = "تنظیمات برنامه پروگرم شده"
$ws.Range("E2").Value = "بررسی مدار Zero Cross"
$ws.Range("F2").Value = "LED"
$ws.Range("G2").Value = "ولتاژ 3.3 ولت"
$ws.Range("H2").Value = "ارتباط سریال"
$ws.Range("I2").Value = "حافظه خارجی"
$ws.Range("J2").Value = "ولتاژ باتری"
$ws.Range("K2").Value = "ولتاز مودم"
$ws.Range("L2").Value = "ارتباط با سیم کارت"
$ws.Range("M2").Value = "کلید"
$ws.Range("N2").Value = "ExcelCheckSumA"
$ws.Range("O2").Value = "ExcelCheckSumB"
$ws.Range("P2").Value = "Version"

# --- Data row (row 4): rewrite the result values ------------------------
$ws.Range("B4").Value = "OK"
$ws.Range("C4").Value = "220v ac:0.05Amp 12v dc:11.72Volt"
$ws.Range("D4").Value = "Chip:'196154487' Domain:'94.139.169.122:8000' Serial:'12345678' Prodoct:'1402-11-23'"
$ws.Range("E4").Value = "OK"
$ws.Range("F4").Value = "Power:OK, RS485:OK, NET1:OK, NET2:OK"
$ws.Range("G4").Value = "OK"
$ws.Range("H4").Value = "OP:OK, RS485:OK"
$ws.Range("I4").Value = "OK"
$ws.Range("J4").Value = "OK"
# J4 used to carry the one-off red-font style; align it with the rest of
# the row (green font) so it collapses back onto the shared style.
$ws.Range("J4").Font.Color = $ws.Range("K4").Font.Color
$ws.Range("K4").Value = "OK"
$ws.Range("L4").Value = "OK"
$ws.Range("M4").Value = "Reset:OK, Factory:OK"
$ws.Range("N4").Value = 1700
$ws.Range("O4").Value = "5F6FCEA73B0B0A331B2988D800CA0DBA"
$ws.Range("P4").Value = "v1.14021121"

# --- Column widths to match the redesigned (wider, text-heavy) table ---
$ws.Columns("C").ColumnWidth = 31.14
$ws.Columns("D").ColumnWidth = 76.43
$ws.Columns("E").ColumnWidth = 33.57
$ws.Columns("F").ColumnWidth = 37.14
$ws.Columns("G").ColumnWidth = 10.71
$ws.Columns("H").ColumnWidth = 17.86
$ws.Columns("I").ColumnWidth = 11
$ws.Columns("J").ColumnWidth = 11
$ws.Columns("K").ColumnWidth = 12.43
$ws.Columns("L").ColumnWidth = 15.14
$ws.Columns("M").ColumnWidth = 21.29
$ws.Columns("N").ColumnWidth = 16
$ws.Columns("O").ColumnWidth = 37.43
$ws.Columns("P").ColumnWidth = 12.57

# --- Selection, as last recorded in the workbook view ------------------
$ws.Range("C6").Select()
